$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental (row 7) -> "false" as TEXT (not boolean).
# Assigning the literal string "false" to .Value auto-converts it to a
# boolean in Excel, so instead we build it as a text formula result and
# paste back as a value, which keeps it typed as text (t="s").
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date (row 8) -> updated timestamp
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Description (row 17) -> new description text
$ws.Range("B17").Value = "Recommended activity types based on recovery status"
